$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the "Status" column (E) values to TRUE for rows 2-5
$ws.Range("E2:E5").Value = $true
